# database/industries/siman/saroum/income/quarterly/dollar_cumulative.xlsx
#
# Quarterly update: the "Overview" sheet is a rolling 10-quarter cumulative
# income-statement window (columns D:M). Each update drops the oldest
# quarter on the left and appends the newest quarter on the right, so every
# data row's contents shift one column to the left and a freshly reported
# quarter ("12 ماهه منتهی به 1401/12") lands in column M. A couple of
# historical publish-date labels (row 9) and one historical figure
# (row 23, discontinued-operations P/L) were also corrected as part of the
# "read_price algorithm" change mentioned in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Row 8: financial-period column headers -- window shifts left by one quarter,
# newest quarter (12 ماهه منتهی به 1401/12) appended in column M.
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# Row 9: publish dates -- shifts left; some mid-window dates were also restated
# (old "1401-10-29" entries corrected to "1401-12-29"/"1402-01-30" variants).
$ws.Range("D9").Value = "1400-10-30 (2)"
$ws.Range("E9").Value = "1401-01-31 (8)"
$ws.Range("F9").Value = "1401-04-30 (2)"
$ws.Range("G9").Value = "1401-08-02 (4)"
$ws.Range("H9").Value = "1401-12-29 (3)"
$ws.Range("I9").Value = "1402-01-30 (9)"
$ws.Range("J9").Value = "1401-04-30"
$ws.Range("K9").Value = "1401-08-02 (2)"
$ws.Range("L9").Value = "1401-12-29 (2)"
$ws.Range("M9").Value = "1402-01-30 (2)"

# Row 11: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D11").Value = 16369
$ws.Range("E11").Value = 19221
$ws.Range("F11").Value = 8069
$ws.Range("G11").Value = 18771
$ws.Range("H11").Value = 24793
$ws.Range("I11").Value = 28449
$ws.Range("J11").Value = 9947
$ws.Range("K11").Value = 18458
$ws.Range("L11").Value = 24983
$ws.Range("M11").Value = 28949

# Row 12: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D12").Value = -7605
$ws.Range("E12").Value = -9341
$ws.Range("F12").Value = -3303
$ws.Range("G12").Value = -6690
$ws.Range("H12").Value = -9846
$ws.Range("I12").Value = -12338
$ws.Range("J12").Value = -4051
$ws.Range("K12").Value = -7438
$ws.Range("L12").Value = -10691
$ws.Range("M12").Value = -13109

# Row 13: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D13").Value = 8763
$ws.Range("E13").Value = 9881
$ws.Range("F13").Value = 4766
$ws.Range("G13").Value = 12081
$ws.Range("H13").Value = 14947
$ws.Range("I13").Value = 16111
$ws.Range("J13").Value = 5896
$ws.Range("K13").Value = 11020
$ws.Range("L13").Value = 14291
$ws.Range("M13").Value = 15841

# Row 14: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D14").Value = -280
$ws.Range("E14").Value = -377
$ws.Range("F14").Value = -312
$ws.Range("G14").Value = -627
$ws.Range("H14").Value = -837
$ws.Range("I14").Value = -1080
$ws.Range("J14").Value = -413
$ws.Range("K14").Value = -665
$ws.Range("L14").Value = -898
$ws.Range("M14").Value = -1404

# Row 16: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D16").Value = 261
$ws.Range("E16").Value = 222
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 59
$ws.Range("H16").Value = 64
$ws.Range("I16").Value = 14
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = 189
$ws.Range("L16").Value = 197
$ws.Range("M16").Value = 1044

# Row 17: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D17").Value = 8744
$ws.Range("E17").Value = 9726
$ws.Range("F17").Value = 4464
$ws.Range("G17").Value = 11513
$ws.Range("H17").Value = 14174
$ws.Range("I17").Value = 15045
$ws.Range("J17").Value = 5487
$ws.Range("K17").Value = 10544
$ws.Range("L17").Value = 13590
$ws.Range("M17").Value = 15480

# Row 19: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D19").Value = 848
$ws.Range("E19").Value = 2621
$ws.Range("F19").Value = 2749
$ws.Range("G19").Value = 2848
$ws.Range("H19").Value = 3151
$ws.Range("I19").Value = 3675
$ws.Range("J19").Value = 398
$ws.Range("K19").Value = 840
$ws.Range("L19").Value = 1127
$ws.Range("M19").Value = 1697

# Row 20: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D20").Value = 9568
$ws.Range("E20").Value = 12323
$ws.Range("F20").Value = 7212
$ws.Range("G20").Value = 14362
$ws.Range("H20").Value = 17325
$ws.Range("I20").Value = 18720
$ws.Range("J20").Value = 5885
$ws.Range("K20").Value = 11368
$ws.Range("L20").Value = 14677
$ws.Range("M20").Value = 17092

# Row 21: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D21").Value = -1604
$ws.Range("E21").Value = -900
$ws.Range("F21").Value = -917
$ws.Range("G21").Value = -2335
$ws.Range("H21").Value = -2586
$ws.Range("I21").Value = -1667
$ws.Range("J21").Value = -853
$ws.Range("K21").Value = -1686
$ws.Range("L21").Value = -2072
$ws.Range("M21").Value = -1383

# Row 22: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D22").Value = 7964
$ws.Range("E22").Value = 11423
$ws.Range("F22").Value = 6296
$ws.Range("G22").Value = 12026
$ws.Range("H22").Value = 14739
$ws.Range("I22").Value = 17052
$ws.Range("J22").Value = 5032
$ws.Range("K22").Value = 9682
$ws.Range("L22").Value = 12605
$ws.Range("M22").Value = 15709

# Row 24: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D24").Value = 7964
$ws.Range("E24").Value = 11423
$ws.Range("F24").Value = 6296
$ws.Range("G24").Value = 12026
$ws.Range("H24").Value = 14739
$ws.Range("I24").Value = 17052
$ws.Range("J24").Value = 5032
$ws.Range("K24").Value = 9682
$ws.Range("L24").Value = 12605
$ws.Range("M24").Value = 15709

# Row 26: numeric data -- shift left by one quarter, append new quarter value in M
$ws.Range("D26").Value = 3182
$ws.Range("E26").Value = 3096
$ws.Range("F26").Value = 3000
$ws.Range("G26").Value = 2825
$ws.Range("H26").Value = 2691
$ws.Range("I26").Value = 2653
$ws.Range("J26").Value = 2382
$ws.Range("K26").Value = 2318
$ws.Range("L26").Value = 3145
$ws.Range("M26").Value = 2834

# Row 18: shift left; several quarters have no reported value ("-")
$ws.Range("D18").Value = -24
$ws.Range("E18").Value = -24
$ws.Range("F18").Value = "-"
$ws.Range("G18").Value = "-"
$ws.Range("H18").Value = "-"
$ws.Range("I18").Value = "-"
$ws.Range("J18").Value = "-"
$ws.Range("K18").Value = -16
$ws.Range("L18").Value = -39
$ws.Range("M18").Value = -85

# Row 23: discontinued-operations line -- mostly "-"; in addition to the shift,
# the 1400/09 figure (now column I) was corrected from "-" to 24 by the updated
# read_price algorithm, and the new quarter (M) reports 16.
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"
$ws.Range("I23").Value = 24
$ws.Range("J23").Value = "-"
$ws.Range("K23").Value = "-"
$ws.Range("L23").Value = "-"
$ws.Range("M23").Value = 16

